# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column G ("K") values recalculated for rows 2-9 and 11 (rows 10 and 12 unchanged)
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 4
$ws.Range("G8").Value = 0
$ws.Range("G9").Value = 1
$ws.Range("G11").Value = 1
